# Referral Portal changes - update the test user suffix used by the
# createUser automation sheet from 20 to 22. The dependent formulas in
# B2 (Test.<n>User) and F2 (thayne.sampson<n>@riomed.com) recalculate
# automatically since they reference A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("createUser")
$ws.Range("A2").Value = 22
